$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$sh = $s.Shapes.Item(1)

# Remove the "Absence du travail en raison de maladie" paragraph (3rd paragraph)
$tr = $sh.TextFrame.TextRange
$para = $tr.Paragraphs(3, 1)
$para.Delete()

# Shrink the shape's height to match the new (shorter) text content
$sh.Height = 212.9534
